# New training stats and stuff
# - Adds two new games: "Ninja" (inserted where "Plunder" used to be, column O)
#   and "Starpilot" (brand new, column Q). The old "Plunder" column data moves
#   to the new column P.
# - Updates a few existing values in the "Leaper" column (L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Give the two new header cells (P1, Q1) the same look as the existing
# bold/bordered header cells (copy format from N1, the "Miner" header).
$ws.Range("N1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("N1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)

# O1 used to be "Plunder" -- it becomes "Ninja". The old Plunder data slides
# over to the brand new P column, and Starpilot is appended in Q.
$ws.Range("O1").Value = "Ninja"
$ws.Range("P1").Value = "Plunder"
$ws.Range("Q1").Value = "Starpilot"

# --- Row 2 (steps) ------------------------------------------------------
$ws.Range("L2").Value = 7946240
$ws.Range("P2").Value = 7946240
$ws.Range("Q2").Value = 7946240

# --- Row 3 (train) -------------------------------------------------------
$ws.Range("L3").Value = 17.41319465637207
$ws.Range("O3").Value = 40.69444274902344
$ws.Range("P3").Value = 4.684027671813965
$ws.Range("Q3").Value = 26.93923568725586

# --- Row 4 (test) -------------------------------------------------------
$ws.Range("L4").Value = 14.91319465637207
$ws.Range("O4").Value = 40.29513931274414
$ws.Range("P4").Value = 4.173611164093018
$ws.Range("Q4").Value = 26.62847137451172
